$d = $word.ActiveDocument

function New-PkgXml([string]$bodyXml) {
    return @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
$bodyXml
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
}

$rPrPlain = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="222222"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr>'
$rPrBold  = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="222222"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr>'

# --- 1. "Using birds with winter and summer counts (i.e migrant birds)" ---
# Split the text so "i.e" sits in its own run wrapped with proofErr spellStart/spellEnd.
$target = $d.Content.Find.Execute("Using birds with winter and summer counts (i.e migrant birds)")
$p1 = $d.Paragraphs.Item(3)
$body1 = "<w:p><w:pPr>$rPrPlain</w:pPr>" +
         "<w:r>$rPrPlain<w:t>Using birds with winter and summer counts (</w:t></w:r>" +
         "<w:proofErr w:type=`"spellStart`"/>" +
         "<w:r>$rPrPlain<w:t>i.e</w:t></w:r>" +
         "<w:proofErr w:type=`"spellEnd`"/>" +
         "<w:r>$rPrPlain<w:t xml:space=`"preserve`"> migrant birds)</w:t></w:r>" +
         "</w:p>"
$p1.Range.InsertXML((New-PkgXml $body1))

# --- 2. "Overview of biodiversity (...)" becomes bold, with spacing after=0 on the paragraph ---
$p2 = $d.Paragraphs.Item(5)
$body2 = "<w:p><w:pPr><w:spacing w:after=`"0`"/>$rPrBold</w:pPr>" +
         "<w:r>$rPrBold<w:t>Overview of biodiversity (importance of biodiversity and current trends).</w:t></w:r>" +
         "</w:p>"
$p2.Range.InsertXML((New-PkgXml $body2))

# --- 3. Insert two new blank paragraphs right after it (one keeps spacing after=0, one without) ---
$p2 = $d.Paragraphs.Item(5)
$insertPoint = $d.Range($p2.Range.End, $p2.Range.End)
$body3 = "<w:p><w:pPr><w:spacing w:after=`"0`"/>$rPrPlain</w:pPr></w:p>" +
         "<w:p><w:pPr>$rPrPlain</w:pPr></w:p>"
$insertPoint.InsertXML((New-PkgXml $body3))

# --- 4. Append two new blank paragraphs at the very end of the document (after "Mention aim.") ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)
$body4 = "<w:p><w:pPr>$rPrPlain</w:pPr></w:p>" +
         "<w:p><w:pPr>$rPrPlain</w:pPr></w:p>"
$endPoint.InsertXML((New-PkgXml $body4))

Write-Host "Final paragraph count:" $d.Paragraphs.Count
